# Auto-generated Excel COM-interop script to apply numeric updates
# to the Anima_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 2375
$ws.Range("J45").Value = 2375
$ws.Range("L45").Value = 7125
$ws.Range("N45").Value = -7509

$ws.Range("H99").Value = 3964.3333
$ws.Range("I99").Value = 553.5454999999999
$ws.Range("J99").Value = 7716.2
$ws.Range("K99").Value = 1660.6365
$ws.Range("L99").Value = 23148.6
$ws.Range("M99").Value = -162.6364999999998
$ws.Range("N99").Value = -26144.6

$ws.Range("H101").Value = 1585.9166
$ws.Range("I101").Value = 1321.2
$ws.Range("J101").Value = 1775
$ws.Range("K101").Value = 3963.6
$ws.Range("L101").Value = 5325
$ws.Range("M101").Value = -2341.6
$ws.Range("N101").Value = -8569

$ws.Range("H103").Value = 63629.625
$ws.Range("I103").Value = 111622
$ws.Range("K103").Value = 334866
$ws.Range("M103").Value = -334280

$ws.Range("H137").Value = 1668150.6
$ws.Range("I137").Value = 2977521.8
$ws.Range("J137").Value = 1678.2273
$ws.Range("K137").Value = 8932565.399999999
$ws.Range("L137").Value = 5034.6819
$ws.Range("M137").Value = -8930015.399999999
$ws.Range("N137").Value = -10134.6819

$ws.Range("H138").Value = 3851.114
$ws.Range("J138").Value = 4009.638
$ws.Range("L138").Value = 12028.914
$ws.Range("N138").Value = -22308.914

$ws.Range("H141").Value = 3853.318
$ws.Range("I141").Value = 1834.5
$ws.Range("J141").Value = 7386.25
$ws.Range("K141").Value = 5503.5
$ws.Range("L141").Value = 22158.75
$ws.Range("M141").Value = -323.5
$ws.Range("N141").Value = -32518.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10350

$ws.Range("H61").Value = 9807209
$ws.Range("I61").Value = 15875043
$ws.Range("K61").Value = 15875043
$ws.Range("M61").Value = -15874831

$ws.Range("H63").Value = 4715.3335
$ws.Range("J63").Value = 4935.091
$ws.Range("L63").Value = 4935.091
$ws.Range("N63").Value = -6307.091

$ws.Range("H66").Value = 4715.3335
$ws.Range("J66").Value = 4935.091
$ws.Range("L66").Value = 24675.455
$ws.Range("N66").Value = -31539.455

$ws.Range("H74").Value = 13515361
$ws.Range("I74").Value = 937.5238000000001
$ws.Range("J74").Value = 31253040
$ws.Range("K74").Value = 937.5238000000001
$ws.Range("L74").Value = 31253040
$ws.Range("M74").Value = -63.52380000000005
$ws.Range("N74").Value = -31254788

$ws.Range("H77").Value = 13515361
$ws.Range("I77").Value = 937.5238000000001
$ws.Range("J77").Value = 31253040
$ws.Range("K77").Value = 4687.619000000001
$ws.Range("L77").Value = 156265200
$ws.Range("M77").Value = -319.6190000000006
$ws.Range("N77").Value = -156273936

$ws.Range("H136").Value = 9807209
$ws.Range("I136").Value = 15875043
$ws.Range("K136").Value = 47625129
$ws.Range("M136").Value = -47622579

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12015
$ws.Range("I82").Value = 5743
$ws.Range("J82").Value = 19332.334
$ws.Range("K82").Value = 5743
$ws.Range("L82").Value = 19332.334
$ws.Range("M82").Value = -5360
$ws.Range("N82").Value = -20098.334

$ws.Range("H85").Value = 12015
$ws.Range("I85").Value = 5743
$ws.Range("J85").Value = 19332.334
$ws.Range("K85").Value = 5743
$ws.Range("L85").Value = 19332.334
$ws.Range("M85").Value = -4417
$ws.Range("N85").Value = -21984.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 300251.5
$ws.Range("I6").Value = 600001
$ws.Range("J6").Value = 502
$ws.Range("K6").Value = 600001
$ws.Range("L6").Value = 502
$ws.Range("M6").Value = -599888
$ws.Range("N6").Value = -728

$ws.Range("H7").Value = 66.5
$ws.Range("I7").Value = 72
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 72
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 41
$ws.Range("N7").Value = -276

$ws.Range("H17").Value = 18000
$ws.Range("J17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15348

$ws.Range("H25").Value = 9013
$ws.Range("J25").Value = 9013
$ws.Range("L25").Value = 9013
$ws.Range("N25").Value = -9361

$ws.Range("H31").Value = 5253.7866
$ws.Range("I31").Value = 1531.0667
$ws.Range("J31").Value = 7146.695
$ws.Range("K31").Value = 1531.0667
$ws.Range("L31").Value = 7146.695
$ws.Range("M31").Value = -1236.0667
$ws.Range("N31").Value = -7736.695

$ws.Range("H34").Value = 5253.7866
$ws.Range("I34").Value = 1531.0667
$ws.Range("J34").Value = 7146.695
$ws.Range("K34").Value = 1531.0667
$ws.Range("L34").Value = 7146.695
$ws.Range("M34").Value = -1329.0667
$ws.Range("N34").Value = -7550.695

$ws.Range("H41").Value = 13579.4
$ws.Range("I41").Value = 4950
$ws.Range("J41").Value = 19332.334
$ws.Range("K41").Value = 4950
$ws.Range("L41").Value = 19332.334
$ws.Range("M41").Value = -4522
$ws.Range("N41").Value = -20188.334

$ws.Range("H50").Value = 16417.8
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 18022.25
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 18022.25
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -19272.25

$ws.Range("H51").Value = 16899.143
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9264

$ws.Range("H59").Value = 17749.25
$ws.Range("J59").Value = 17749.25
$ws.Range("L59").Value = 17749.25
$ws.Range("N59").Value = -20039.25

$ws.Range("H60").Value = 13334.083
$ws.Range("I60").Value = 1520
$ws.Range("J60").Value = 21772.715
$ws.Range("K60").Value = 1520
$ws.Range("L60").Value = 21772.715
$ws.Range("M60").Value = -1009
$ws.Range("N60").Value = -22794.715

$ws.Range("H61").Value = 16899.143
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9652

$ws.Range("H68").Value = 23272.182
$ws.Range("J68").Value = 23272.182
$ws.Range("L68").Value = 23272.182
$ws.Range("N68").Value = -24770.182

$ws.Range("H71").Value = 23272.182
$ws.Range("J71").Value = 23272.182
$ws.Range("L71").Value = 69816.546
$ws.Range("N71").Value = -77304.546

$ws.Range("H74").Value = 19699.5
$ws.Range("J74").Value = 19699.5
$ws.Range("L74").Value = 19699.5
$ws.Range("N74").Value = -21447.5

$ws.Range("H77").Value = 19699.5
$ws.Range("J77").Value = 19699.5
$ws.Range("L77").Value = 59098.5
$ws.Range("N77").Value = -67834.5

$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 757.1667
$ws.Range("I92").Value = 616.6667
$ws.Range("J92").Value = 897.6667
$ws.Range("K92").Value = 1850.0001
$ws.Range("L92").Value = 2693.0001
$ws.Range("M92").Value = -602.0001
$ws.Range("N92").Value = -5189.0001

$ws.Range("H129").Value = 705946.6
$ws.Range("I129").Value = 467.66666
$ws.Range("J129").Value = 979035.25
$ws.Range("K129").Value = 1402.99998
$ws.Range("L129").Value = 2937105.75
$ws.Range("M129").Value = 3597.00002
$ws.Range("N129").Value = -2947105.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 55564124
$ws.Range("I132").Value = 83343680
$ws.Range("J132").Value = 5004.1665
$ws.Range("K132").Value = 250031040
$ws.Range("L132").Value = 15012.4995
$ws.Range("M132").Value = -250028510
$ws.Range("N132").Value = -20072.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 36665.4
$ws.Range("J70").Value = 36665.4
$ws.Range("L70").Value = 36665.4
$ws.Range("N70").Value = -37205.4

$ws.Range("H73").Value = 36665.4
$ws.Range("J73").Value = 36665.4
$ws.Range("L73").Value = 36665.4
$ws.Range("N73").Value = -38537.4

$ws.Range("H74").Value = 21428.572
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002

$ws.Range("H77").Value = 21428.572
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008
